# Fix formatting issues introduced when scraping floating point numbers:
#  1) A handful of "Razon social"/"Nombre Fantasia" entries used a comma as a
#     generic separator between multiple names; replace those commas with
#     periods (these are exact, literal text fixes).
#  2) The "Importe" (amount) column stored numbers using Spanish/"es-AR"
#     formatting (thousands separator "." and decimal separator ","), e.g.
#     "7.440,00". These need to become plain dot-decimal numbers stored as
#     text, e.g. "7440.00".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Literal text fixes (comma -> period) -------------------------------
$ws.Cells.Replace("SCHAB DARIO, PEROTTI XAVIER, BENINCA MATIAS S.H.", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH") | Out-Null
$ws.Cells.Replace("BOFFELLI, MARIA INES", "BOFFELLI. MARIA INES") | Out-Null
$ws.Cells.Replace("MARSICO GUILLERMO MIGUEL, MARSICO JUAN EDUARDO", "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO") | Out-Null
$ws.Cells.Replace("PARPAGNOLI, PEDRO RICARDO", "PARPAGNOLI. PEDRO RICARDO") | Out-Null
$ws.Cells.Replace("ALBIZZATTI, PABLO MARTIN Y FULINI, SERGIO RUBEN", "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN") | Out-Null

# --- 2) Reformat "Importe" column (H) numbers -------------------------------
# These cells are stored as text (shared strings), not real numbers, so we
# must be careful to keep them as text after rewriting. Writing a plain
# numeric-looking string back via .Value/.Formula makes Excel coerce the
# cell into a genuine number (losing the literal ".00" and the text type),
# so instead we build a text-formula ("=""7440.00""") for every cell and
# then convert the whole column to literal values in one shot, which keeps
# the cells as plain text without touching any cell styles/number formats.

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1
$col = 8 # column H = "Importe"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $orig = $cell.Text
    if ($orig -ne "") {
        $newVal = $orig.Replace(".", "").Replace(",", ".")
        $cell.Formula = '="' + $newVal + '"'
    }
}

$colRange = $ws.Range($ws.Cells.Item(2, $col), $ws.Cells.Item($lastRow, $col))
$colRange.Copy()
$colRange.PasteSpecial(-4163) # xlPasteValues
$excel.CutCopyMode = 0
